$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.189.44'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.94%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.658.72'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.57%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5146'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.03%  '
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2581'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.97%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06437'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.98'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07823'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.74%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.300'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.30%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.652.95'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.885.38'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5552'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8065'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.93%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.26'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.208.53'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.91%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '211.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.430'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.977'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("E24").Value = '  +0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.06'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.756'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1166'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.976'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.80'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05210'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.254'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.367'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.65%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.220'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.570'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.759'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.63%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9320'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.32%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.373'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.174.40'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +12.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5699'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.80%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01593'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8445'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.58%  '
$ws.Range("E42").Value = '  +0.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.674'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.62'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.795.73'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.66%  '
$ws.Range("E46").Value = '  +1.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4536'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.90'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.07%  '
$ws.Range("E49").Value = '  -0.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.894'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.87%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05063'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.17%  '
